$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update facility_or_program_type text (column B) for rows whose alphabetical
# position changed after the dataset re-sort
$ws.Range("B3").Value = "City-financed homeless set-aside units units financed jointly by HPD and HDC under Housing New York"
$ws.Range("B4").Value = "City-financed homeless set-aside units- HDC set- aside units financed under Housing New York"
$ws.Range("B5").Value = "City-financed homeless set-aside units- HPD set-aside units financed under Housing New York"
$ws.Range("B24").Value = "Private rental market apartment with a rental subsidy- disaggregated by the type of such subsidy (CFHEPS)"
$ws.Range("B25").Value = "Private rental market apartment with a rental subsidy- disaggregated by the type of such subsidy (EHV)"
$ws.Range("B26").Value = "Private rental market apartment with a rental subsidy- disaggregated by the type of such subsidy (FHEPS)"
$ws.Range("B27").Value = "Private rental market apartment with a rental subsidy- disaggregated by the type of such subsidy (SOTA)"
$ws.Range("B29").Value = "Rapid re-housing funded by the United States department of housing and urban development"
$ws.Range("B30").Value = "Residential drug treatment and detoxification"
$ws.Range("B31").Value = "Section 8"
$ws.Range("B32").Value = "Section 8 voucher housing- HPD- project-based"
$ws.Range("B33").Value = "Section 8 voucher housing- HPD- tenant-based"
$ws.Range("B34").Value = "Section 8 voucher housing- NYC Housing Authority (NYCHA)- project-based"
$ws.Range("B35").Value = "Section 8 voucher housing- NYCHA- project-based"
$ws.Range("B36").Value = "Section 8 voucher housing- NYCHA- tenant-based"
$ws.Range("B37").Value = "Section 8 voucher housing- New York state homes and community renewal- tenant-based"
$ws.Range("B38").Value = "Section 8 voucher housing- New York state homes and community renewal-project-based"
$ws.Range("B39").Value = "Settings with higher levels of medical care- inpatient hospitalization"
$ws.Range("B40").Value = "Settings with higher levels of medical care- long-term care facilities"
$ws.Range("B41").Value = "Settings with higher levels of medical care- medical rehabilitation centers"
$ws.Range("B42").Value = "Settings with higher levels of medical care- medical respite care"
$ws.Range("B43").Value = "Shared Living (Not friends or relatives)"
$ws.Range("B44").Value = "Subsidized Apartment (NYCHA, Mitchell Lama, Etc.)"
$ws.Range("B45").Value = "Supportive housing"
$ws.Range("B46").Value = "Transitional housing operated by or under contract or similar agreement with DHS, DYCD, HPD, United States Department of Housing and Urban Development (HUD) or HRA"
$ws.Range("B47").Value = "Transitional housing operated by or under contract or similar agreement with DHS, DYCD, HPD, United States department of housing and urban development or HRA"
$ws.Range("B48").Value = "Unknown or unable to validate"
$ws.Range("B49").Value = "Youth detention center/Correctional facility"
$ws.Range("B50").Value = "private rental market apartment with a rental subsidy- disaggregated by the type of such subsidy"
$ws.Range("B51").Value = "settings with higher levels of medical care- inpatient hospitalization"
$ws.Range("B52").Value = "settings with higher levels of medical care- long-term care facilities"
$ws.Range("B53").Value = "settings with higher levels of medical care- medical rehabilitation centers"
$ws.Range("B54").Value = "settings with higher levels of medical care- medical respite care"

# Update the matching n (count) column so each label keeps its correct count
$ws.Range("C24").Value = 7.0
$ws.Range("C29").Value = 28.0
$ws.Range("C31").Value = 14.0
$ws.Range("C32").Value = 28.0
$ws.Range("C34").Value = 21.0
$ws.Range("C35").Value = 7.0
$ws.Range("C37").Value = 28.0
$ws.Range("C38").Value = 28.0
$ws.Range("C39").Value = 21.0
$ws.Range("C40").Value = 21.0
$ws.Range("C42").Value = 21.0
$ws.Range("C43").Value = 7.0
$ws.Range("C45").Value = 28.0
$ws.Range("C46").Value = 21.0
$ws.Range("C47").Value = 7.0
$ws.Range("C48").Value = 28.0
$ws.Range("C50").Value = 21.0
$ws.Range("C51").Value = 7.0
$ws.Range("C53").Value = 7.0
